$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2274052478134111
$ws.Range("C2").Value = 0.5014577259475219
$ws.Range("J2").Value = 0.01457725947521866
$ws.Range("P2").Value = 0.1690962099125364
$ws.Range("S2").Value = 0.08746355685131195
$ws.Range("B3").Value = 0.01104972375690608
$ws.Range("C3").Value = 0.02762430939226519
$ws.Range("J3").Value = 0.03867403314917127
$ws.Range("P3").Value = 0.7182320441988951
$ws.Range("S3").Value = 0.2044198895027624
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6382978723404256
$ws.Range("S4").Value = 0.3404255319148936
$ws.Range("B6").Value = 0.04958677685950413
$ws.Range("D6").Value = 0.01652892561983471
$ws.Range("F6").Value = 0.04958677685950413
$ws.Range("J6").Value = 0.3181818181818182
$ws.Range("O6").Value = 0.02479338842975207
$ws.Range("Q6").Value = 0.128099173553719
$ws.Range("R6").Value = 0.06611570247933884
$ws.Range("S6").Value = 0.3471074380165289
$ws.Range("B7").Value = 0.09734513274336283
$ws.Range("D7").Value = 0.004424778761061947
$ws.Range("E7").Value = 0.004424778761061947
$ws.Range("F7").Value = 0.05309734513274336
$ws.Range("J7").Value = 0.1150442477876106
$ws.Range("O7").Value = 0.02212389380530973
$ws.Range("Q7").Value = 0.1814159292035398
$ws.Range("R7").Value = 0.08849557522123894
$ws.Range("S7").Value = 0.4336283185840708
$ws.Range("B8").Value = 0.1092150170648464
$ws.Range("D8").Value = 0.01535836177474403
$ws.Range("F8").Value = 0.07167235494880546
$ws.Range("J8").Value = 0.1040955631399317
$ws.Range("O8").Value = 0.01877133105802048
$ws.Range("Q8").Value = 0.1535836177474403
$ws.Range("R8").Value = 0.07337883959044368
$ws.Range("S8").Value = 0.4539249146757679
$ws.Range("B9").Value = 0.09793814432989691
$ws.Range("D9").Value = 0.02061855670103093
$ws.Range("F9").Value = 0.05154639175257732
$ws.Range("J9").Value = 0.1030927835051546
$ws.Range("O9").Value = 0.02577319587628866
$ws.Range("Q9").Value = 0.1649484536082474
$ws.Range("R9").Value = 0.09278350515463918
$ws.Range("S9").Value = 0.4432989690721649
$ws.Range("B10").Value = 0.09881697981906751
$ws.Range("D10").Value = 0.02157272094641614
$ws.Range("E10").Value = 0.00208768267223382
$ws.Range("F10").Value = 0.06123869171885873
$ws.Range("J10").Value = 0.1210855949895616
$ws.Range("O10").Value = 0.01600556715379262
$ws.Range("Q10").Value = 0.1899791231732777
$ws.Range("R10").Value = 0.08837856645789841
$ws.Range("S10").Value = 0.4008350730688935
$ws.Range("G11").Value = 0.1573033707865168
$ws.Range("J11").Value = 0.08707865168539326
$ws.Range("K11").Value = 0.2134831460674157
$ws.Range("L11").Value = 0.5224719101123596
$ws.Range("S11").Value = 0.01966292134831461
$ws.Range("G12").Value = 0.7643979057591623
$ws.Range("J12").Value = 0.1884816753926702
$ws.Range("K12").Value = 0.01047120418848168
$ws.Range("L12").Value = 0.02094240837696335
$ws.Range("S12").Value = 0.01570680628272251
$ws.Range("G13").Value = 0.62
$ws.Range("J13").Value = 0.28
$ws.Range("S13").Value = 0.1
$ws.Range("F15").Value = 0.03448275862068965
$ws.Range("H15").Value = 0.228448275862069
$ws.Range("I15").Value = 0.06465517241379311
$ws.Range("J15").Value = 0.2887931034482759
$ws.Range("K15").Value = 0.04741379310344827
$ws.Range("M15").Value = 0.01724137931034483
$ws.Range("O15").Value = 0.04741379310344827
$ws.Range("S15").Value = 0.271551724137931
$ws.Range("F16").Value = 0.01851851851851852
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("I16").Value = 0.09722222222222222
$ws.Range("J16").Value = 0.4305555555555556
$ws.Range("K16").Value = 0.09722222222222222
$ws.Range("M16").Value = 0.01851851851851852
$ws.Range("O16").Value = 0.05092592592592592
$ws.Range("S16").Value = 0.1203703703703704
$ws.Range("F17").Value = 0.01068376068376068
$ws.Range("H17").Value = 0.2243589743589744
$ws.Range("J17").Value = 0.4123931623931624
$ws.Range("K17").Value = 0.1004273504273504
$ws.Range("M17").Value = 0.01495726495726496
$ws.Range("O17").Value = 0.0576923076923077
$ws.Range("S17").Value = 0.1239316239316239
$ws.Range("F18").Value = 0.02678571428571428
$ws.Range("H18").Value = 0.2455357142857143
$ws.Range("I18").Value = 0.1026785714285714
$ws.Range("J18").Value = 0.3526785714285715
$ws.Range("K18").Value = 0.1026785714285714
$ws.Range("M18").Value = 0.02232142857142857
$ws.Range("O18").Value = 0.05803571428571429
$ws.Range("S18").Value = 0.08928571428571429
$ws.Range("F19").Value = 0.02213541666666667
$ws.Range("H19").Value = 0.2252604166666667
$ws.Range("I19").Value = 0.07096354166666667
$ws.Range("J19").Value = 0.3736979166666667
$ws.Range("K19").Value = 0.1145833333333333
$ws.Range("M19").Value = 0.02213541666666667
$ws.Range("N19").Value = 0.0006510416666666666
$ws.Range("O19").Value = 0.05729166666666666
$ws.Range("S19").Value = 0.11328125
